$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 1.98
$ws.Range("G3").Value = 2.2
$ws.Range("I3").Value = 4.2
$ws.Range("K3").Value = 4.3
$ws.Range("P3").Value = 2.12
$ws.Range("Q3").Value = 1.61
